$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.2086021505376344
$ws.Cells.Item(2, 3).Value = 0.5419354838709678
$ws.Cells.Item(2, 10).Value = 0.01720430107526882
$ws.Cells.Item(2, 16).Value = 0.1268817204301075
$ws.Cells.Item(2, 19).Value = 0.1053763440860215
$ws.Cells.Item(3, 2).Value = 0.0078125
$ws.Cells.Item(3, 3).Value = 0.00390625
$ws.Cells.Item(3, 10).Value = 0.01953125
$ws.Cells.Item(3, 16).Value = 0.75390625
$ws.Cells.Item(3, 19).Value = 0.21484375
$ws.Cells.Item(4, 10).Value = 0.06578947368421052
$ws.Cells.Item(4, 15).Value = 0.0131578947368421
$ws.Cells.Item(4, 16).Value = 0.6842105263157895
$ws.Cells.Item(4, 19).Value = 0.2368421052631579
$ws.Cells.Item(6, 2).Value = 0.06201550387596899
$ws.Cells.Item(6, 4).Value = 0.01937984496124031
$ws.Cells.Item(6, 5).Value = 0.003875968992248062
$ws.Cells.Item(6, 6).Value = 0.07364341085271318
$ws.Cells.Item(6, 10).Value = 0.3837209302325582
$ws.Cells.Item(6, 15).Value = 0.01162790697674419
$ws.Cells.Item(6, 17).Value = 0.1317829457364341
$ws.Cells.Item(6, 18).Value = 0.04263565891472868
$ws.Cells.Item(6, 19).Value = 0.2713178294573643
$ws.Cells.Item(7, 2).Value = 0.1273584905660377
$ws.Cells.Item(7, 4).Value = 0.02358490566037736
$ws.Cells.Item(7, 6).Value = 0.05188679245283019
$ws.Cells.Item(7, 10).Value = 0.1320754716981132
$ws.Cells.Item(7, 15).Value = 0.0330188679245283
$ws.Cells.Item(7, 17).Value = 0.1367924528301887
$ws.Cells.Item(7, 18).Value = 0.1037735849056604
$ws.Cells.Item(7, 19).Value = 0.3915094339622642
$ws.Cells.Item(8, 2).Value = 0.1577777777777778
$ws.Cells.Item(8, 4).Value = 0.02666666666666667
$ws.Cells.Item(8, 5).Value = 0.002222222222222222
$ws.Cells.Item(8, 6).Value = 0.05111111111111111
$ws.Cells.Item(8, 10).Value = 0.1088888888888889
$ws.Cells.Item(8, 15).Value = 0.02444444444444445
$ws.Cells.Item(8, 17).Value = 0.1644444444444444
$ws.Cells.Item(8, 18).Value = 0.09333333333333334
$ws.Cells.Item(8, 19).Value = 0.3711111111111111
$ws.Cells.Item(9, 2).Value = 0.1261261261261261
$ws.Cells.Item(9, 4).Value = 0.02252252252252252
$ws.Cells.Item(9, 5).Value = 0.004504504504504504
$ws.Cells.Item(9, 6).Value = 0.03153153153153153
$ws.Cells.Item(9, 10).Value = 0.1351351351351351
$ws.Cells.Item(9, 15).Value = 0.02252252252252252
$ws.Cells.Item(9, 17).Value = 0.2072072072072072
$ws.Cells.Item(9, 18).Value = 0.07657657657657657
$ws.Cells.Item(9, 19).Value = 0.3738738738738739
$ws.Cells.Item(10, 2).Value = 0.1509695290858726
$ws.Cells.Item(10, 4).Value = 0.03393351800554017
$ws.Cells.Item(10, 5).Value = 0.002077562326869806
$ws.Cells.Item(10, 6).Value = 0.05678670360110803
$ws.Cells.Item(10, 10).Value = 0.1191135734072022
$ws.Cells.Item(10, 15).Value = 0.0131578947368421
$ws.Cells.Item(10, 17).Value = 0.195983379501385
$ws.Cells.Item(10, 18).Value = 0.08518005540166206
$ws.Cells.Item(10, 19).Value = 0.342797783933518
$ws.Cells.Item(11, 7).Value = 0.1238390092879257
$ws.Cells.Item(11, 10).Value = 0.0804953560371517
$ws.Cells.Item(11, 11).Value = 0.1764705882352941
$ws.Cells.Item(11, 12).Value = 0.6006191950464397
$ws.Cells.Item(11, 19).Value = 0.01857585139318885
$ws.Cells.Item(12, 7).Value = 0.7828282828282829
$ws.Cells.Item(12, 10).Value = 0.1515151515151515
$ws.Cells.Item(12, 11).Value = 0.0101010101010101
$ws.Cells.Item(12, 12).Value = 0.0202020202020202
$ws.Cells.Item(12, 19).Value = 0.03535353535353535
$ws.Cells.Item(13, 7).Value = 0.6097560975609756
$ws.Cells.Item(13, 10).Value = 0.3414634146341464
$ws.Cells.Item(13, 19).Value = 0.04878048780487805
$ws.Cells.Item(15, 6).Value = 0.0339622641509434
$ws.Cells.Item(15, 8).Value = 0.1018867924528302
$ws.Cells.Item(15, 9).Value = 0.07169811320754717
$ws.Cells.Item(15, 10).Value = 0.3962264150943396
$ws.Cells.Item(15, 11).Value = 0.07924528301886792
$ws.Cells.Item(15, 13).Value = 0.003773584905660377
$ws.Cells.Item(15, 15).Value = 0.06792452830188679
$ws.Cells.Item(15, 19).Value = 0.2452830188679245
$ws.Cells.Item(16, 6).Value = 0.02013422818791946
$ws.Cells.Item(16, 8).Value = 0.174496644295302
$ws.Cells.Item(16, 9).Value = 0.09060402684563758
$ws.Cells.Item(16, 10).Value = 0.3456375838926175
$ws.Cells.Item(16, 11).Value = 0.1073825503355705
$ws.Cells.Item(16, 13).Value = 0.02348993288590604
$ws.Cells.Item(16, 14).Value = 0.006711409395973154
$ws.Cells.Item(16, 15).Value = 0.07718120805369127
$ws.Cells.Item(16, 19).Value = 0.1543624161073825
$ws.Cells.Item(17, 6).Value = 0.05353319057815846
$ws.Cells.Item(17, 8).Value = 0.158458244111349
$ws.Cells.Item(17, 9).Value = 0.07922912205567452
$ws.Cells.Item(17, 10).Value = 0.4025695931477516
$ws.Cells.Item(17, 11).Value = 0.08565310492505353
$ws.Cells.Item(17, 13).Value = 0.01498929336188437
$ws.Cells.Item(17, 14).Value = 0.002141327623126338
$ws.Cells.Item(17, 15).Value = 0.08565310492505353
$ws.Cells.Item(17, 19).Value = 0.1177730192719486
$ws.Cells.Item(18, 6).Value = 0.04265402843601896
$ws.Cells.Item(18, 8).Value = 0.1137440758293839
$ws.Cells.Item(18, 9).Value = 0.1421800947867299
$ws.Cells.Item(18, 10).Value = 0.4454976303317535
$ws.Cells.Item(18, 11).Value = 0.08530805687203792
$ws.Cells.Item(18, 13).Value = 0.01895734597156398
$ws.Cells.Item(18, 15).Value = 0.04265402843601896
$ws.Cells.Item(18, 19).Value = 0.1090047393364929
$ws.Cells.Item(19, 6).Value = 0.02633504023408925
$ws.Cells.Item(19, 8).Value = 0.1982443306510607
$ws.Cells.Item(19, 9).Value = 0.08266276517922458
$ws.Cells.Item(19, 10).Value = 0.3672275054864667
$ws.Cells.Item(19, 11).Value = 0.1126554498902707
$ws.Cells.Item(19, 13).Value = 0.0182882223847842
$ws.Cells.Item(19, 15).Value = 0.07461594732991954
$ws.Cells.Item(19, 19).Value = 0.1199707388441843
